# Auto-generated edit script applying the Raiden_Profits.xlsx diff
# Sets updated LeveProfit / market-price values across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1864.2222
$ws.Range("J17").Value = 1864.2222
$ws.Range("L17").Value = 5592.6666
$ws.Range("N17").Value = -5928.6666
$ws.Range("H33").Value = 425.97675
$ws.Range("I33").Value = 425.97675
$ws.Range("K33").Value = 425.97675
$ws.Range("M33").Value = -196.97675
$ws.Range("H55").Value = 503.77777
$ws.Range("I55").Value = 166.8
$ws.Range("K55").Value = 166.8
$ws.Range("M55").Value = 47.19999999999999
$ws.Range("H87").Value = 17600
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
$ws.Range("H90").Value = 17600
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
$ws.Range("H113").Value = 10059.462
$ws.Range("I113").Value = 10343.637
$ws.Range("K113").Value = 10343.637
$ws.Range("M113").Value = -7089.637000000001
$ws.Range("H125").Value = 2664.6667
$ws.Range("J125").Value = 2491.4
$ws.Range("L125").Value = 22422.6
$ws.Range("N125").Value = -27342.6
$ws.Range("H132").Value = 1840
$ws.Range("I132").Value = 1480.1666
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 4440.4998
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -1910.4998
$ws.Range("N132").Value = -17057
$ws.Range("H137").Value = 7699.9
$ws.Range("I137").Value = 3692.1428
$ws.Range("J137").Value = 9857.923000000001
$ws.Range("K137").Value = 11076.4284
$ws.Range("L137").Value = 29573.769
$ws.Range("M137").Value = -8526.428400000001
$ws.Range("N137").Value = -34673.769
$ws.Range("H138").Value = 9754.531999999999
$ws.Range("J138").Value = 9847.604499999999
$ws.Range("L138").Value = 29542.8135
$ws.Range("N138").Value = -39822.8135

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23770.062
$ws.Range("I32").Value = 15423.654
$ws.Range("J32").Value = 59937.832
$ws.Range("K32").Value = 15423.654
$ws.Range("L32").Value = 59937.832
$ws.Range("M32").Value = -15136.654
$ws.Range("N32").Value = -60511.832
$ws.Range("H46").Value = 18866
$ws.Range("J46").Value = 23333
$ws.Range("L46").Value = 23333
$ws.Range("N46").Value = -23971
$ws.Range("H61").Value = 6599.5654
$ws.Range("I61").Value = 5988
$ws.Range("J61").Value = 7160.1665
$ws.Range("K61").Value = 5988
$ws.Range("L61").Value = 7160.1665
$ws.Range("M61").Value = -5776
$ws.Range("N61").Value = -7584.1665
$ws.Range("H63").Value = 2676.1667
$ws.Range("I63").Value = 2632.4055
$ws.Range("K63").Value = 2632.4055
$ws.Range("M63").Value = -1946.4055
$ws.Range("H66").Value = 2676.1667
$ws.Range("I66").Value = 2632.4055
$ws.Range("K66").Value = 13162.0275
$ws.Range("M66").Value = -9730.0275
$ws.Range("H74").Value = 1888.625
$ws.Range("I74").Value = 1856
$ws.Range("K74").Value = 1856
$ws.Range("M74").Value = -982
$ws.Range("H77").Value = 1888.625
$ws.Range("I77").Value = 1856
$ws.Range("K77").Value = 9280
$ws.Range("M77").Value = -4912
$ws.Range("H110").Value = 3288.4443
$ws.Range("I110").Value = 3288.4443
$ws.Range("K110").Value = 3288.4443
$ws.Range("M110").Value = -1243.4443
$ws.Range("H122").Value = 2057.5518
$ws.Range("I122").Value = 1999.5652
$ws.Range("K122").Value = 5998.6956
$ws.Range("M122").Value = -3548.6956
$ws.Range("H136").Value = 6599.5654
$ws.Range("I136").Value = 5988
$ws.Range("J136").Value = 7160.1665
$ws.Range("K136").Value = 17964
$ws.Range("L136").Value = 21480.4995
$ws.Range("M136").Value = -15414
$ws.Range("N136").Value = -26580.4995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1684.75
$ws.Range("J20").Value = 2225.6
$ws.Range("L20").Value = 2225.6
$ws.Range("N20").Value = -2719.6
$ws.Range("H105").Value = 6522.7827
$ws.Range("I105").Value = 5772
$ws.Range("K105").Value = 5772
$ws.Range("M105").Value = -4025
$ws.Range("H107").Value = 1797.7273
$ws.Range("I107").Value = 1797.7273
$ws.Range("K107").Value = 1797.7273
$ws.Range("M107").Value = 122.2727
$ws.Range("H134").Value = 5095.28
$ws.Range("I134").Value = 3707.182
$ws.Range("J134").Value = 6185.9287
$ws.Range("K134").Value = 11121.546
$ws.Range("L134").Value = 18557.7861
$ws.Range("M134").Value = -8586.545999999998
$ws.Range("N134").Value = -23627.7861

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3313.3333
$ws.Range("I31").Value = 3313.3333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3313.3333
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3018.3333
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 3313.3333
$ws.Range("I34").Value = 3313.3333
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3313.3333
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3111.3333
$ws.Range("N34").ClearContents()
$ws.Range("H41").Value = 6511.8
$ws.Range("I41").Value = 6511.8
$ws.Range("K41").Value = 6511.8
$ws.Range("M41").Value = -6083.8
$ws.Range("H62").Value = 8389.6
$ws.Range("I62").Value = 8781.375
$ws.Range("K62").Value = 8781.375
$ws.Range("M62").Value = -8157.375
$ws.Range("H65").Value = 8389.6
$ws.Range("I65").Value = 8781.375
$ws.Range("K65").Value = 43906.875
$ws.Range("M65").Value = -40786.875
$ws.Range("H107").Value = 424.0909
$ws.Range("I107").Value = 236.5
$ws.Range("J107").Value = 2300
$ws.Range("K107").Value = 236.5
$ws.Range("L107").Value = 2300
$ws.Range("M107").Value = 1683.5
$ws.Range("N107").Value = -6140
$ws.Range("H141").Value = 308415.25
$ws.Range("J141").Value = 331907.53
$ws.Range("L141").Value = 331907.53
$ws.Range("N141").Value = -342267.53

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5552
$ws.Range("J39").Value = 7000
$ws.Range("L39").Value = 21000
$ws.Range("N39").Value = -21588
$ws.Range("H109").Value = 1838.25
$ws.Range("I109").Value = 1838.25
$ws.Range("K109").Value = 5514.75
$ws.Range("M109").Value = -4474.75
$ws.Range("H131").Value = 27146.111
$ws.Range("J131").Value = 7014.5
$ws.Range("L131").Value = 21043.5
$ws.Range("N131").Value = -31123.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1098.8667
$ws.Range("I97").Value = 1076.2222
$ws.Range("J97").Value = 1132.8334
$ws.Range("K97").Value = 1076.2222
$ws.Range("L97").Value = 1132.8334
$ws.Range("M97").Value = -580.2221999999999
$ws.Range("N97").Value = -2124.8334
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2241.3333
$ws.Range("I46").Value = 1249.6666
$ws.Range("J46").Value = 3233
$ws.Range("K46").Value = 1249.6666
$ws.Range("L46").Value = 3233
$ws.Range("M46").Value = -1061.6666
$ws.Range("N46").Value = -3609
$ws.Range("H55").Value = 1023.75
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H95").Value = 49975
$ws.Range("J95").Value = 49975
$ws.Range("L95").Value = 49975
$ws.Range("N95").Value = -55467
$ws.Range("H105").Value = 26910
$ws.Range("J105").Value = 26910
$ws.Range("L105").Value = 26910
$ws.Range("N105").Value = -33898

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 37499.2
$ws.Range("I31").Value = 25999.2
$ws.Range("J31").Value = 48999.2
$ws.Range("K31").Value = 25999.2
$ws.Range("L31").Value = 48999.2
$ws.Range("M31").Value = -25651.2
$ws.Range("N31").Value = -49695.2
$ws.Range("H62").Value = 30829.9
$ws.Range("I62").Value = 12966.333
$ws.Range("J62").Value = 38485.715
$ws.Range("K62").Value = 12966.333
$ws.Range("L62").Value = 38485.715
$ws.Range("M62").Value = -12342.333
$ws.Range("N62").Value = -39733.715
$ws.Range("H65").Value = 30829.9
$ws.Range("I65").Value = 12966.333
$ws.Range("J65").Value = 38485.715
$ws.Range("K65").Value = 64831.665
$ws.Range("L65").Value = 192428.575
$ws.Range("M65").Value = -61711.665
$ws.Range("N65").Value = -198668.575
$ws.Range("H126").Value = 2288.6191
$ws.Range("I126").Value = 2470.7778
$ws.Range("K126").Value = 7412.3334
$ws.Range("M126").Value = -4942.3334
$ws.Range("H132").Value = 3305.2917
$ws.Range("I132").Value = 3230.0952
$ws.Range("J132").Value = 3831.6667
$ws.Range("K132").Value = 9690.285600000001
$ws.Range("L132").Value = 11495.0001
$ws.Range("M132").Value = -7160.285600000001
$ws.Range("N132").Value = -16555.0001
